$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Delete the "codigo_funcionario" (Codigo do funcionario da creche) row ---
# Row 30 used to describe a field linking Aluno to Funcionario via a
# "Codigo do funcionario da creche" column; that field + its FK note are removed.
$ws.Range("A30:D30").ClearContents()

# --- Add a new field to the Funcionario table (row 54) ---
# Replaces the removed FK with a CNPJ-based foreign key to Creche.
$ws.Range("A54").Value2 = "CNPJ"
$ws.Range("B54").Value2 = "varchar(14)"
$ws.Range("C54").Value2 = "CNPJ da creche"
$ws.Range("D54").Value2 = "Not Null; FK"

# --- Formatting tweak on B25 (touched while editing the selection) ---
$ws.Range("B25").Font.Bold = $true
$ws.Range("B25").Font.Bold = $false

# --- Update the view: scroll position + active selection ---
$win = $excel.ActiveWindow
$ws.Range("A13").Select()
$ws.Range("B25").Select()
$win.ScrollRow = 13
